# Refresh the cryptos list: update Price (column D) and Volume(1h) (column E)
# text values for rows 2-51 on the active worksheet to match the latest scrape.
#
# D-column entries are plain text cells (e.g. "27.654.34", "208.22") even
# though some of them look like numbers. Typing a numeric-looking value
# straight into a cell would make Excel auto-convert it to a real number, so
# each D value is entered with a leading apostrophe (forces text entry) and
# the cell's style is then reset to "Normal" so no stray number-format
# styling sticks around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = "<price text>"; E = "<+/-X.XX%>" }
$updates = @{
    2  = @{ D = "27.649.24";  E = "+0.33%" }
    3  = @{ D = "1.589.26";   E = "-0.34%" }
    4  = @{ E = "+0.20%" }
    5  = @{ D = "208.03";     E = "+0.24%" }
    6  = @{ D = "0.500";      E = "-0.16%" }
    7  = @{ E = "+0.21%" }
    8  = @{ D = "22.25";      E = "-0.07%" }
    9  = @{ E = "-0.38%" }
    10 = @{ E = "+0.19%" }
    11 = @{ D = "0.0866";     E = "-0.58%" }
    12 = @{ D = "1.814.95";   E = "-0.40%" }
    13 = @{ D = "1.579.28";   E = "-1.17%" }
    14 = @{ E = "-0.86%" }
    15 = @{ E = "-2.26%" }
    16 = @{ D = "27.665.39";  E = "+0.48%" }
    17 = @{ E = "-0.38%" }
    18 = @{ D = "216.92";     E = "-0.10%" }
    19 = @{ E = "-0.05%" }
    20 = @{ E = "-0.83%" }
    21 = @{ E = "+0.13%" }
    22 = @{ E = "-1.22%" }
    23 = @{ D = "9.76";       E = "+0.70%" }
    24 = @{ D = "1.99";       E = "-0.49%" }
    25 = @{ D = "153.22";     E = "-1.36%" }
    26 = @{ D = "6.98";       E = "+4.11%" }
    27 = @{ E = "+0.21%" }
    28 = @{ D = "15.05";      E = "+0.32%" }
    29 = @{ E = "-0.66%" }
    30 = @{ E = "-0.09%" }
    31 = @{ D = "0.0473";     E = "+1.21%" }
    32 = @{ E = "-2.33%" }
    33 = @{ D = "1.374.21";   E = "+1.59%" }
    34 = @{ E = "+0.67%" }
    35 = @{ D = "1.53";       E = "-0.32%" }
    36 = @{ E = "+0.60%" }
    37 = @{ D = "2.31";       E = "+0.02%" }
    38 = @{ E = "+1.75%" }
    39 = @{ D = "0.533";      E = "-1.02%" }
    40 = @{ D = "0.824";      E = "+1.06%" }
    41 = @{ E = "+0.18%" }
    42 = @{ D = "0.976";      E = "+2.33%" }
    43 = @{ E = "+0.40%" }
    44 = @{ E = "+4.10%" }
    45 = @{ D = "1.76";       E = "+0.71%" }
    46 = @{ E = "-1.70%" }
    47 = @{ D = "1.727.37";   E = "-0.20%" }
    48 = @{ D = "85.84";      E = "-1.67%" }
    49 = @{ E = "+1.19%" }
    50 = @{ D = "0.0962";     E = "-0.62%" }
    51 = @{ D = "0.0494";     E = "-0.50%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.Value = "'" + $vals["D"]
        $cell.Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = "  " + $vals["E"] + "  "
    }
}
